$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts old N/O/P -> O/P/Q),
# mirroring Excel's native "Insert Column" behaviour (new column inherits
# the formatting of the column immediately to its left).
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 9.75

# Make "Repayment schedule" the active sheet (was "Transactions"), and
# update the current selection on it.
$ws.Activate()
[void]$ws.Range("R11").Select()
